$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1366576124457447
$ws.Range("C2").Value = 0.775372048323969
$ws.Range("D2").Value = 1.844561190472676
$ws.Range("E2").Value = 1.358146233095934
$ws.Range("F2").Value = 1.364436723157444

# Row 3
$ws.Range("B3").Value = 0.100644678693309
$ws.Range("C3").Value = 0.5870089618847463
$ws.Range("D3").Value = 0.9743935215277069
$ws.Range("E3").Value = 0.9871137328229745
$ws.Range("F3").Value = 0.9917406180962824

# Row 4
$ws.Range("B4").Value = 0.09787366134911048
$ws.Range("C4").Value = 0.6307949273637553
$ws.Range("D4").Value = 1.168272060271323
$ws.Range("E4").Value = 1.080866347089835
$ws.Range("F4").Value = 1.087354403430002

# Row 5
$ws.Range("B5").Value = 0.08154868731045999
$ws.Range("C5").Value = 0.713998770293546
$ws.Range("D5").Value = 1.430771530282745
$ws.Range("E5").Value = 1.196148623826799
$ws.Range("F5").Value = 1.208190537716507
$ws.Range("G5").Value = 41

# Row 6
$ws.Range("B6").Value = 0.1341520404658785
$ws.Range("C6").Value = 0.7691896598824904
$ws.Range("D6").Value = 1.579558871789514
$ws.Range("E6").Value = 1.256805025367703
$ws.Range("F6").Value = 1.270281138917405
$ws.Range("G6").Value = 31

# Row 7
$ws.Range("B7").Value = 0.1268598637344239
$ws.Range("C7").Value = 0.7813345819487484
$ws.Range("D7").Value = 1.592716119073123
$ws.Range("E7").Value = 1.262028573001865
$ws.Range("F7").Value = 1.277861748157917
$ws.Range("G7").Value = 29

# Row 8
$ws.Range("B8").Value = 0.09675479195317427
$ws.Range("C8").Value = 0.8108202300352334
$ws.Range("D8").Value = 1.726960037846746
$ws.Range("E8").Value = 1.314138515471922
$ws.Range("F8").Value = 1.335537356534567
$ws.Range("G8").Value = 27

# Row 9
$ws.Range("B9").Value = 0.1406368155206148
$ws.Range("C9").Value = 0.8950670472051415
$ws.Range("D9").Value = 2.222009981461377
$ws.Range("E9").Value = 1.490640795584697
$ws.Range("F9").Value = 1.524656502007553
$ws.Range("G9").Value = 19

# Row 10
$ws.Range("B10").Value = -0.3312026136842777
$ws.Range("C10").Value = 0.8013776111497077
$ws.Range("D10").Value = 1.293401115164856
$ws.Range("E10").Value = 1.137277941035021
$ws.Range("F10").Value = 1.136360671227658
$ws.Range("G10").Value = 12

# Row 11
$ws.Range("B11").Value = 0.4641315885571715
$ws.Range("C11").Value = 0.7592743174351484
$ws.Range("D11").Value = 1.498317365315948
$ws.Range("E11").Value = 1.224057745907418
$ws.Range("F11").Value = 1.266342782296397
